# Regenerate save_data: column G ("K" - strike count) recomputed, and row 9 H/I/J corrected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 1
    4 = 0
    5 = 1
    6 = 1
    7 = 2
    8 = 1
    9 = 2
    10 = 3
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 3
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 1
    31 = 1
    32 = 3
    33 = 1
    34 = 3
    35 = 0
    36 = 1
    37 = 1
    38 = 2
    39 = 1
    40 = 1
    41 = 3
    42 = 2
    43 = 1
    44 = 0
    45 = 2
    46 = 0
    47 = 2
    48 = 1
    49 = 0
    50 = 0
    51 = 0
    52 = 2
    53 = 1
    54 = 1
    55 = 1
    56 = 0
    57 = 0
    58 = 1
    59 = 1
    60 = 1
    61 = 2
    62 = 0
    63 = 0
    64 = 0
    65 = 3
    66 = 1
    67 = 2
    68 = 2
    69 = 1
    70 = 1
    71 = 1
    72 = 2
    73 = 1
    74 = 2
    75 = 2
    76 = 2
    77 = 1
    78 = 1
    79 = 2
    80 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

# Row 9 also had IP/I0/IF (H/I/J) recomputed during the regen
$ws.Cells.Item(9, 8).Value = 3
$ws.Cells.Item(9, 9).Value = 6
$ws.Cells.Item(9, 10).Value = 8
